# 2022FSAdates.xlsx - "Add files via upload"
# The workbook's physical worksheet file xl/worksheets/sheet2.xml is exposed
# under the Excel sheet tab named "Sheet1" (sheet name <-> file name are
# swapped in this workbook), so we explicitly grab that worksheet by name
# rather than relying on ActiveSheet/Index ordering.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# --- Data edits (Table1 columns: D=Cancelled, E=SplitGender) -------------

# Row 15 (20220410 / U17-20 / F): SplitGender 1 -> 0
$ws.Range("E15").Value = 0

# Row 16 (20220410 / U15 / S): Cancelled 0 -> 1
$ws.Range("D16").Value = 1

# Row 18 (20220410 / U17-20 / S): Cancelled 0 -> 1
$ws.Range("D18").Value = 1

# Row 20 (20220410 / U17-20 / E): Cancelled 0 -> 1
$ws.Range("D20").Value = 1

# Row 21 (20220410 / OB / S): Cancelled 0 -> 1
$ws.Range("D21").Value = 1

# --- View state: selection changed to D15:L21 (active cell D15) ----------
# (topLeftCell scroll position A10 is an Excel-window-chrome setting that
# this host does not expose/persist, but the selection itself is applied.)

$win = $excel.ActiveWindow
$win.ScrollRow = 10
$win.ScrollColumn = 1
$ws.Range("D15:L21").Select()
